$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/01_kitipi"
$ws.Range("B2").Value = "pngimages/01_gift.png"
$ws.Range("C2").Value = "trainingimages/16_kotapi"
$ws.Range("D2").Value = "pngimages/16_icecream.png"
$ws.Range("E2").Value = -0.5
$ws.Range("F2").Value = 0.5

# Row 3
$ws.Range("A3").Value = "trainingimages/05_titopo"
$ws.Range("B3").Value = "pngimages/05_megaphone.png"
$ws.Range("C3").Value = "trainingimages/26_kapako"
$ws.Range("D3").Value = "pngimages/26_pineapple.png"

# Row 4
$ws.Range("A4").Value = "trainingimages/15_kopota"
$ws.Range("B4").Value = "pngimages/15_barrel.png"
$ws.Range("C4").Value = "trainingimages/23_patoko"
$ws.Range("D4").Value = "pngimages/23_lemon.png"
